# Auto-generated Excel COM-interop script to apply the diff
# Update latest output (run 204)

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet (rows 2-5, columns A-F) ---
$wsSchedule.Range("A2").Value = 46081.04166666666
$wsSchedule.Range("B2").Value = 46081.22916666666
$wsSchedule.Range("C2").Value = 4.5
$wsSchedule.Range("D2").Value = 17.01
$wsSchedule.Range("E2").Value = 491.6404544999999
$wsSchedule.Range("F2").Value = 28.90302495590829
$wsSchedule.Range("A3").Value = 46081.29166666666
$wsSchedule.Range("B3").Value = 46081.77083333334
$wsSchedule.Range("C3").Value = 11.5
$wsSchedule.Range("D3").Value = 43.47
$wsSchedule.Range("E3").Value = 208.37638575
$wsSchedule.Range("F3").Value = 4.793567650103521
$wsSchedule.Range("A4").Value = 46082.0625
$wsSchedule.Range("B4").Value = 46082.25
$wsSchedule.Range("C4").Value = 4.5
$wsSchedule.Range("D4").Value = 17.01
$wsSchedule.Range("E4").Value = 217.55314425
$wsSchedule.Range("F4").Value = 12.78972041446208
$wsSchedule.Range("A5").Value = 46082.29166666666
$wsSchedule.Range("B5").Value = 46082.60416666666
$wsSchedule.Range("C5").Value = 7.5
$wsSchedule.Range("D5").Value = 28.35
$wsSchedule.Range("E5").Value = 148.40035275
$wsSchedule.Range("F5").Value = 5.234580343915344

# --- Detailed sheet ---
$wsDetailed.Range("E4").Value = "ON"
$wsDetailed.Range("E5").Value = "ON"
$wsDetailed.Range("E6").Value = "ON"
$wsDetailed.Range("E7").Value = "ON"
$wsDetailed.Range("E8").Value = "ON"
$wsDetailed.Range("E9").Value = "ON"
$wsDetailed.Range("E10").Value = "ON"
$wsDetailed.Range("E11").Value = "ON"
$wsDetailed.Range("E12").Value = "ON"
$wsDetailed.Range("E16").Value = "ON"
$wsDetailed.Range("B37").Value = 37.89
$wsDetailed.Range("B38").Value = 50.3379
$wsDetailed.Range("B39").Value = 57.36
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("E39").Value = "OFF"
$wsDetailed.Range("B40").Value = 58.87904
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("E40").Value = "OFF"
$wsDetailed.Range("B41").Value = 62.82274
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 61.16608
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("E42").Value = "OFF"
$wsDetailed.Range("B43").Value = 58.91223
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("E43").Value = "OFF"
$wsDetailed.Range("B44").Value = 57.36
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("E44").Value = "OFF"
$wsDetailed.Range("B45").Value = 59.55755
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("E45").Value = "OFF"
$wsDetailed.Range("B46").Value = 57.32
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("E46").Value = "OFF"
$wsDetailed.Range("B47").Value = 57.31
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "OFF"
$wsDetailed.Range("B48").Value = 47.04295
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("E48").Value = "OFF"
$wsDetailed.Range("B49").Value = 51.03101
$wsDetailed.Range("E49").Value = "OFF"
$wsDetailed.Range("B51").Value = 51.74779
$wsDetailed.Range("E52").Value = "OFF"
$wsDetailed.Range("B53").Value = 25.70682
$wsDetailed.Range("B54").Value = 23.73563
$wsDetailed.Range("B55").Value = 23.76635
$wsDetailed.Range("B56").Value = 23.92999
$wsDetailed.Range("B57").Value = 24.37851
$wsDetailed.Range("B58").Value = 24.53589
$wsDetailed.Range("B59").Value = 26.14065
$wsDetailed.Range("B60").Value = 26.13224
$wsDetailed.Range("E60").Value = "ON"
$wsDetailed.Range("B61").Value = 24.80535
$wsDetailed.Range("E61").Value = "ON"
$wsDetailed.Range("B62").Value = 26.75644
$wsDetailed.Range("B63").Value = 40.54
$wsDetailed.Range("B64").Value = 10.22666
$wsDetailed.Range("E64").Value = "ON"
$wsDetailed.Range("B67").Value = 0.78944
$wsDetailed.Range("B69").Value = 0.7
$wsDetailed.Range("B70").Value = 2.68644
$wsDetailed.Range("B72").Value = 2.77295
$wsDetailed.Range("B73").Value = 22.07
$wsDetailed.Range("B74").Value = 22.07
$wsDetailed.Range("B75").Value = 22.07
$wsDetailed.Range("B76").Value = 22.07
$wsDetailed.Range("B77").Value = 22.07
$wsDetailed.Range("B78").Value = 22.07
$wsDetailed.Range("B79").Value = 27.11344
$wsDetailed.Range("E79").Value = "OFF"
$wsDetailed.Range("B80").Value = 37.89
$wsDetailed.Range("E80").Value = "OFF"
$wsDetailed.Range("B83").Value = 37.89
$wsDetailed.Range("B84").Value = 51.00322
$wsDetailed.Range("B85").Value = 53.26765
$wsDetailed.Range("B86").Value = 57.31
$wsDetailed.Range("B87").Value = 57.36
$wsDetailed.Range("B90").Value = 57.06007
$wsDetailed.Range("B91").Value = 57.06
$wsDetailed.Range("B92").Value = 52.31016
$wsDetailed.Range("B93").Value = 51.50835
$wsDetailed.Range("B94").Value = 50.15574
$wsDetailed.Range("B95").Value = 50.10404
$wsDetailed.Range("B96").Value = 37.89
$wsDetailed.Range("B97").Value = 49.48127

Write-Output "Applied all changes"
